# Automatische test-sync: 2025-06-27 22:28:50
# Adds a new log row (row 7) to the "Logs" sheet describing test mail #1
# ("Wanneer zijn jullie open?") that was answered at 2025-06-27 22:28:38,
# extends the conditional formatting ranges to cover the new row, and
# bumps the "Openingstijden / Locatie" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet -----------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 7

# The new entry is another "Wanneer zijn jullie open?" test mail, identical
# in content to the previous occurrences already logged in rows 4-6 - only
# the timestamp differs. Re-use the existing values so formatting/line
# breaks line up exactly with the sibling rows.
$logs.Cells.Item($newRow, 1).Value = $logs.Cells.Item(4, 1).Value()
$logs.Cells.Item($newRow, 2).Value = $logs.Cells.Item(4, 2).Value()
$logs.Cells.Item($newRow, 3).Value = $logs.Cells.Item(4, 3).Value()
$logs.Cells.Item($newRow, 4).Value = $logs.Cells.Item(4, 4).Value()
$logs.Cells.Item($newRow, 5).Value = $logs.Cells.Item(4, 5).Value()
$logs.Cells.Item($newRow, 6).Value = "2025-06-27 22:28:38"
$logs.Cells.Item($newRow, 7).Value = $logs.Cells.Item(4, 7).Value()
$logs.Cells.Item($newRow, 8).Value = $logs.Cells.Item(4, 8).Value()
$logs.Cells.Item($newRow, 9).Value = $logs.Cells.Item(4, 9).Value()

# Writing multi-line text into a brand new row makes the engine apply a
# custom row height (mirroring Excel's auto-fit-on-entry behaviour). Reset
# it back to the sheet's standard auto height so row 7 stays plain, like
# every other data row.
$logs.Rows.Item($newRow).AutoFit()

# Grow every conditional-formatting block that watched D2:D6 / G2:G6 /
# H2:H6 / I2:I6 so it also covers the freshly added row 7.
$ranges = @("D2:D6", "G2:G6", "H2:H6", "I2:I6")
foreach ($addr in $ranges) {
    $rng = $logs.Range($addr)
    $fcs = $rng.FormatConditions
    $count = $fcs.Count()
    $col = $addr.Substring(0, 1)
    $newAddr = $col + "2:" + $col + "7"
    for ($i = 1; $i -le $count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newAddr))
    }
}

# --- Dashboard sheet --------------------------------------------------------
# One more "Openingstijden / Locatie" reply was logged, so its running
# total goes from 3 to 4.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 4
